# [Fonds de solidarite] Add 2020-11-10 data
#
# The source sheet stores every data cell (columns C "nombre_aides" and
# D "montant_total") as literal TEXT, not as numbers, even though the
# content looks numeric. A plain `Range.Value = "154"` assignment lets
# the usual Excel type inference kick in and silently turns the cell
# into a Number, which would change the cell's stored type relative to
# the original workbook. To keep the cells as Text we:
#   1. force the cell's NumberFormat to Text ("@") so the assignment is
#      not re-interpreted as a number,
#   2. assign the new value,
#   3. reset the cell Style back to "Normal" so we don't leave a stray
#      text-only style behind on the cell (keeps formatting identical
#      to before the edit; only the content changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @("C35", "154"),
    @("D35", "607791.72"),

    @("C36", "762"),
    @("D36", "3269577.29"),

    @("C37", "358"),
    @("D37", "2662305.18"),

    @("C43", "233"),
    @("D43", "806226.81"),

    @("C44", "99"),
    @("D44", "477429.00"),

    @("C45", "31"),
    @("D45", "308095.14"),

    @("C47", "95"),
    @("D47", "313908.00"),

    @("C62", "1157"),
    @("D62", "3687757.20"),

    @("C64", "3188"),
    @("D64", "19048085.47"),

    @("C65", "1105"),
    @("D65", "8081490.96"),

    @("C110", "409"),
    @("D110", "1277943.68"),

    @("C111", "1657"),
    @("D111", "6572255.16"),

    @("C112", "661"),
    @("D112", "4103463.08"),

    @("C113", "215"),
    @("D113", "2241621.73"),

    @("C114", "74"),
    @("D114", "931840.00"),

    @("C115", "10"),
    @("D115", "254339.00")
)

foreach ($edit in $edits) {
    $cellRef = $edit[0]
    $newValue = $edit[1]

    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $newValue
    $range.Style = "Normal"
}
